$d = $word.ActiveDocument

# Locate the run that holds the page break right after
# "...plusieurs fonctionnalités, " and replace it with the new sentence
# "la plus logique est l'affichage avec son écran de 5,7 pouces avec
# résolution haute définition. " (the page break itself is removed).
$r = $d.Content
$found = $r.Find.Execute("fonctionnalités, ^m", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Narrow the range down to just the page-break character (the last
    # character matched by the search above) and drop it.
    $r.Start = $r.End - 1
    $r.Delete()

    # Insert the replacement text in place of the removed page break.
    $r.InsertAfter("la plus logique est l’affichage avec son écran de 5,7 pouces avec résolution haute définition. ")
}
